# Apply the "Flugbahnen angepasst. Neue Pizzasymbole. Erweiterungen im XML
# erfasst. Kommentare erfasst." edit to the Entwicklung backlog sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New rows 24-27: newly captured comments / backlog entries
# (written first so the shared-string table grows in the same order as the
# original author's edit: new "Essen" comments before the reworded
# "Restart Button" text)
$ws.Cells.Item(24, 1).Value = 23
$ws.Cells.Item(24, 2).Value = "Wenn ein Insekt über das Essen fliegt soll ein Geräusch entstehen, z.B. Rülpsen"
$ws.Cells.Item(24, 3).Value = "open"
$ws.Cells.Item(24, 4).Value = "med"
$ws.Rows.Item(24).RowHeight = 30

# --- Row 14: "Restart Button" -> "Restart Button, ESC-Button.", tested -> ok
$ws.Cells.Item(14, 2).Value = "Restart Button, ESC-Button."
$ws.Cells.Item(14, 6).Value = "ok"

# --- Row 15: now also marked as tested -> ok
$ws.Cells.Item(15, 6).Value = "ok"

# --- Row 22/23: flight-path backlog items closed & tested
$ws.Cells.Item(22, 3).Value = "closed"
$ws.Cells.Item(22, 6).Value = "ok"

$ws.Cells.Item(23, 3).Value = "closed"
$ws.Cells.Item(23, 6).Value = "ok"

# --- New row 25
$ws.Cells.Item(25, 1).Value = 24
$ws.Cells.Item(25, 2).Value = "Insektensummgeräusch"
$ws.Cells.Item(25, 3).Value = "open"
$ws.Cells.Item(25, 4).Value = "med"

# --- New row 26
$ws.Cells.Item(26, 1).Value = 25
$ws.Cells.Item(26, 2).Value = "Wenn auf das Essen geklickt wird, soll GameOver erscheinen"
$ws.Cells.Item(26, 3).Value = "open"
$ws.Cells.Item(26, 4).Value = "low"
$ws.Rows.Item(26).RowHeight = 30

# --- New row 27
$ws.Cells.Item(27, 1).Value = 26
$ws.Cells.Item(27, 2).Value = "Essen soll nach jedem Insektenüberflug, stückweise kleiner werden."
$ws.Cells.Item(27, 3).Value = "open"
$ws.Cells.Item(27, 4).Value = "med"
$ws.Rows.Item(27).RowHeight = 30

# --- Re-apply the autofilter over the grown table range A1:I26
$ws.AutoFilterMode = $false
$ws.Range("A1:I26").AutoFilter()

# --- Keep the _FilterDatabase defined name in sync with the new range
foreach ($n in $wb.Names) {
    if ($n.Name -eq "Entwicklung!_FilterDatabase") {
        $n.RefersTo = "=Entwicklung!`$A`$1:`$I`$26"
    }
}

# --- Update the saved selection to match the new bottom of the list
$ws.Range("E27").Select()
